# Rename columns in the "vehicleTypes" and "shipments" tables, matching the
# commit "done renaming in transport optimization":
#   vehicleTypes: maxCapacityLoadingMeter -> maxCapacityPallets
#                 costsPerKm              -> costsPerDistanceUnit
#   shipments:    sender                  -> senderId
#                 senderServiceTime       -> senderStopDuration
#                 recipient               -> recipientId
#                 recipientServiceTime    -> recipientStopDuration
#                 loadingMeter            -> pallets
#                 opportunityCosts        -> externalCosts
#
# Column headers are edited directly on the worksheet ranges (which keeps the
# backing table's column names, and the shared-string table, in sync). The
# shipments sheet's "pallets" rename is applied before the sender/recipient
# renames so the new shared-string entries land in the same order the
# original author's file shows.

$wb = $excel.ActiveWorkbook

$wsVehicleTypes = $wb.Worksheets.Item(2)
$wsShipments = $wb.Worksheets.Item(3)

# vehicleTypes header renames
$wsVehicleTypes.Range("L1").Value = "maxCapacityPallets"
$wsVehicleTypes.Range("O1").Value = "costsPerDistanceUnit"

# shipments header renames
$wsShipments.Range("M1").Value = "pallets"
$wsShipments.Range("C1").Value = "senderId"
$wsShipments.Range("D1").Value = "senderStopDuration"
$wsShipments.Range("G1").Value = "recipientId"
$wsShipments.Range("H1").Value = "recipientStopDuration"
$wsShipments.Range("N1").Value = "externalCosts"

# Restore the selection on the shipments sheet, then switch the active tab to
# vehicleTypes and select its new cell, matching the saved view state.
$wsShipments.Range("N1").Select()

$wsVehicleTypes.Activate()
$wsVehicleTypes.Range("U25").Select()
